# issue #5: property land done
# Normalises the "土地" (land) sheet onto the standard export schema
# (name/area/share_portion/owner/register_date/register_reason/
#  acquire_value/property_category/category/date/legislator_name/
#  legislator_id/source_file/index), drops the stray garbled second
# data row, and cleans up a handful of OCR/whitespace artefacts that
# were left over in the other property sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 土地 (land)
# ---------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

# Header row: rename the existing columns and clean their text, then
# extend the header with the shared pipeline columns.
$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

# New trailing header cells - copy the header format (bold + border)
# from H1 onto I1:O1 before writing their text.
$land.Range("H1").Copy()
$land.Range("I1:O1").PasteSpecial(-4122)
$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# Data row 2: fix whitespace/formatting glitches in the scraped text.
$land.Range("B2").Value = "新北市永和區永福段08230000地號"
$land.Range("F2").Value = "80年06月13日"
$land.Range("G2").Value = "第一次登記"

# New trailing data cells - copy the data format from H2 onto I2:O2
# first so they pick up the plain (unbordered) data style.
$land.Range("H2").Copy()
$land.Range("I2:O2").PasteSpecial(-4122)
$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
# "2011-12-18" must stay a literal text value, not get auto-parsed
# into a date serial - enter it with a leading apostrophe to force
# text, then re-paste a plain sibling's format to drop the resulting
# quote-prefix style.
$land.Range("K2").Value = "'2011-12-18"
$land.Range("J2").Copy()
$land.Range("K2").PasteSpecial(-4122)
$land.Range("L2").Value = "林淑芬"
$land.Range("M2").Value = 1337
$land.Range("N2").Value = "tmp7b501"
$land.Range("O2").Value = 15

# Drop the old, garbled third row entirely (it was a mis-parsed
# duplicate of the header, not real data).
$land.Rows.Item(3).Delete()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Sheet 2: 建物 (building) - text clean-up, plus the running "index"
# column (A) shifts down by 4 because the land sheet above now only
# keeps a single data row instead of two.
# ---------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Range("A2").Value = 20
$building.Range("B2").Value = "新北市永和區永福段02017000建號"
$building.Range("F2").Value = "80年06月130"
$building.Range("G2").Value = "第一次登記"
$building.Range("A3").Value = 22
$building.Range("A4").Value = 23
$building.Range("D4").Value = "面積（平方公尺）"
$building.Range("E4").Value = "權利範圍(持分）"

# ---------------------------------------------------------------
# Sheet 3: 汽車 (car)
# ---------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Range("A2").Value = 34
$car.Range("B2").Value = "toyotarav4rod"
$car.Range("E2").Value = "100年03月01曰"
$car.Range("F2").Value = "nrry=*=r貝賣"

# ---------------------------------------------------------------
# Sheet 4: 現金 (cash)
# ---------------------------------------------------------------
$cash = $wb.Worksheets.Item("現金")
$cash.Range("D1").Value = "新臺幣總額或折合新臺幣總額"
$cash.Range("A2").Value = 44
$cash.Range("C2").Value = "邱若山"
$cash.Range("A3").Value = 45

# D3 needs to stay a genuine text value ("1527000") rather than be
# coerced back into a number - enter it with a leading apostrophe to
# force text, then re-paste the plain sibling cell's format so the
# cell doesn't keep a quote-prefix style.
$cash.Range("D3").Value = "'1527000"
$cash.Range("D2").Copy()
$cash.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$cash.Range("A4").Value = 46

# ---------------------------------------------------------------
# Sheet 5: 存款 (bank deposits)
# ---------------------------------------------------------------
$deposit = $wb.Worksheets.Item("存款")
$deposit.Range("A2").Value = 50
$deposit.Range("A3").Value = 51
$deposit.Range("B3").Value = "遠東國際商業銀行永和分行"
$deposit.Range("A4").Value = 52
$deposit.Range("B4").Value = "遠東國際商業銀行永和分行"

Write-Host "land sheet normalised"
